$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.622.36"
$ws.Range("E2").Value = "  +1.31%  "

# Row 3
$ws.Range("D3").Value = "3.395.35"
$ws.Range("E3").Value = "  +0.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'578.02"
$ws.Range("E5").Value = "  +1.07%  "

# Row 6
$ws.Range("D6").Value = "'137.44"
$ws.Range("E6").Value = "  +2.06%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.395.94"
$ws.Range("E8").Value = "  +1.00%  "

# Row 9
$ws.Range("D9").Value = "'0.473"
$ws.Range("E9").Value = "  -0.77%  "

# Row 10
$ws.Range("D10").Value = "'7.51"
$ws.Range("E10").Value = "  -1.27%  "

# Row 11
$ws.Range("E11").Value = "  +2.72%  "

# Row 12
$ws.Range("D12").Value = "'0.391"
$ws.Range("E12").Value = "  +0.31%  "

# Row 13
$ws.Range("D13").Value = "3.976.07"
$ws.Range("E13").Value = "  +1.06%  "

# Row 14
$ws.Range("E14").Value = "  +2.29%  "

# Row 15
$ws.Range("E15").Value = "  +2.94%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.394.82"
$ws.Range("E16").Value = "  +0.97%  "

# Row 17
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'25.93"
$ws.Range("E17").Value = "  +3.32%  "

# Row 18
$ws.Range("D18").Value = "61.712.51"
$ws.Range("E18").Value = "  +1.28%  "

# Row 19
$ws.Range("D19").Value = "'14.20"
$ws.Range("E19").Value = "  +2.17%  "

# Row 20
$ws.Range("D20").Value = "'5.89"
$ws.Range("E20").Value = "  +2.27%  "

# Row 21
$ws.Range("D21").Value = "'9.44"
$ws.Range("E21").Value = "  +0.02%  "

# Row 22
$ws.Range("D22").Value = "'376.52"
$ws.Range("E22").Value = "  +1.10%  "

# Row 23
$ws.Range("D23").Value = "'0.558"
$ws.Range("E23").Value = "  -2.67%  "

# Row 24
$ws.Range("D24").Value = "3.534.14"
$ws.Range("E24").Value = "  +1.06%  "

# Row 25
$ws.Range("D25").Value = "'0.0000127"
$ws.Range("E25").Value = "  +9.07%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
$ws.Range("D27").Value = "'71.25"
$ws.Range("E27").Value = "  +0.66%  "

# Row 28
$ws.Range("D28").Value = "'1.67"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("D29").Value = "'7.54"
$ws.Range("E29").Value = "  -1.76%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.161"
$ws.Range("E31").Value = "  +4.25%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.25"
$ws.Range("E32").Value = "  +1.83%  "

# Row 33
$ws.Range("E33").Value = "  +1.97%  "

# Row 35
$ws.Range("D35").Value = "'23.44"
$ws.Range("E35").Value = "  +0.42%  "

# Row 36
$ws.Range("D36").Value = "'5.33"
$ws.Range("E36").Value = "  -3.96%  "

# Row 37
$ws.Range("D37").Value = "'1.55"
$ws.Range("E37").Value = "  +0.32%  "

# Row 38
$ws.Range("D38").Value = "'6.84"
$ws.Range("E38").Value = "  -1.03%  "

# Row 39
$ws.Range("D39").Value = "'165.64"
$ws.Range("E39").Value = "  +1.93%  "

# Row 40
$ws.Range("D40").Value = "'0.0781"
$ws.Range("E40").Value = "  -0.29%  "

# Row 41
$ws.Range("B41").Value = "ONDO"
$ws.Range("C41").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D41").Value = "'1.24"
$ws.Range("E41").Value = "  +2.87%  "

# Row 42
$ws.Range("D42").Value = "'0.782"
$ws.Range("E42").Value = "  +3.35%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.74"
$ws.Range("E43").Value = "  +8.76%  "

# Row 44
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("E45").Value = "  +9.61%  "

# Row 46
$ws.Range("D46").Value = "'4.42"
$ws.Range("E46").Value = "  +0.56%  "

# Row 47
$ws.Range("D47").Value = "'41.50"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("E48").Value = "  -1.37%  "

# Row 49
$ws.Range("D49").Value = "'22.73"
$ws.Range("E49").Value = "  -2.52%  "

# Row 50
$ws.Range("D50").Value = "2.337.34"
$ws.Range("E50").Value = "  +5.39%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0262"
$ws.Range("E51").Value = "  +2.02%  "
